$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.014.66"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "2.061.45"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'249.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").Value = "'0.673"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.75%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'54.84"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.53%  "
$ws.Range("D9").Value = "'61.09"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("D10").Value = "'0.382"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("D11").Value = "'0.0798"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.97%  "
$ws.Range("E12").Value = "  +5.81%  "
$ws.Range("E13").Value = "  +2.71%  "
$ws.Range("D14").Value = "2.363.25"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").Value = "'0.817"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").Value = "'5.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.74%  "
$ws.Range("D17").Value = "2.057.59"
$ws.Range("E17").Value = "  -1.97%  "
$ws.Range("D18").Value = "36.971.92"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").Value = "0.0₃0943"
$ws.Range("E19").Value = "  +13.00%  "
$ws.Range("D20").Value = "'73.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "'14.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.07%  "
$ws.Range("D22").Value = "'5.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.57%  "
$ws.Range("D23").Value = "'237.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("D25").Value = "'2.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.77%  "
$ws.Range("D26").Value = "'170.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("D28").Value = "'20.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.31%  "
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("D31").Value = "'4.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.62%  "
$ws.Range("D32").Value = "'0.0629"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").Value = "'1.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.90%  "
$ws.Range("E34").Value = "  +7.08%  "
$ws.Range("D35").Value = "'0.0886"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  -6.81%  "
$ws.Range("E38").Value = "  -5.20%  "
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("E40").Value = "  +22.63%  "
$ws.Range("D41").Value = "'17.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.32%  "
$ws.Range("D42").Value = "'0.0225"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("E43").Value = "  -2.18%  "
$ws.Range("D44").Value = "'96.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.24%  "
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").Value = "'4.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +36.31%  "
$ws.Range("D47").Value = "'13.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -52.29%  "
$ws.Range("D48").Value = "'2.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.33%  "
$ws.Range("D49").Value = "1.297.92"
$ws.Range("E49").Value = "  -3.10%  "
$ws.Range("E50").Value = "  +1.10%  "
$ws.Range("D51").Value = "'4.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.15%  "
